$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values look numeric, so force text formatting to avoid
# Excel auto-converting them to numbers (preserving original string formatting).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.465.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.993"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0870"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.789.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.563.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.442.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "225.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.107"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.445.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.540"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.701.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.67"
$ws.Range("D48").Style = "Normal"

# Coin name, link, and volume columns are already treated as text by Excel.
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  -1.88%  "
